$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column L, matching the style of the other header cells (B1:K1 -> style index 1)
$ws.Range("L1").Value = "Umrank"
$ws.Range("K1").Copy()
$ws.Range("L1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# New numeric data for L2:L12
$values = @(7.5, 9, 7.5, 10, 2.5, 4, 2.5, 11, 1, 5, 6)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 12).Value = $values[$i]
}
